$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = '65.689.13'
$ws.Range('E2').Value = '  -1.02%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '3.450.61'
$ws.Range('E3').Value = '  -3.80%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  -0.04%  '

# Row 5 - BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.67'
$ws.Range('E5').Value = '  -1.44%  '

# Row 6 - Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.68'
$ws.Range('E6').Value = '  -7.73%  '

# Row 7 - LidoStakedEther
$ws.Range('D7').Value = '3.448.75'
$ws.Range('E7').Value = '  -3.83%  '

# Row 8 - USDC
$ws.Range('E8').Value = '  -0.12%  '

# Row 9 - XRP
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.492'
$ws.Range('E9').Value = '  +0.46%  '

# Row 10 - Toncoin
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.57'
$ws.Range('E10').Value = '  -4.43%  '

# Row 11 - Dogecoin
$ws.Range('E11').Value = '  -9.85%  '

# Row 12 - Cardano
$ws.Range('E12').Value = '  -7.94%  '

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range('D13').Value = '4.029.66'
$ws.Range('E13').Value = '  -4.02%  '

# Row 14 - ShibaInu
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000182'
$ws.Range('E14').Value = '  -11.23%  '

# Row 15 - rotate->WrappedEther
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.489.19'
$ws.Range('E15').Value = '  -2.86%  '

# Row 16 - rotate->Avalanche
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.54'
$ws.Range('E16').Value = '  -10.06%  '

# Row 17 - rotate->WrappedBTC
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '65.584.91'
$ws.Range('E17').Value = '  -1.25%  '

# Row 18 - TRON
$ws.Range('E18').Value = '  -2.21%  '

# Row 19 - Uniswap
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.91'
$ws.Range('E19').Value = '  -10.53%  '

# Row 20 - Polkadot
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.78'
$ws.Range('E20').Value = '  -8.73%  '

# Row 21 - Chainlink
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.77'
$ws.Range('E21').Value = '  -7.27%  '

# Row 22 - BitcoinCash
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '395.57'
$ws.Range('E22').Value = '  -6.37%  '

# Row 23 - Polygon
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.548'
$ws.Range('E23').Value = '  -10.22%  '

# Row 24 - Litecoin
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.32'
$ws.Range('E24').Value = '  -6.17%  '

# Row 25 - Dai
$ws.Range('E25').Value = '  -0.14%  '

# Row 26 - WrappedeETH
$ws.Range('D26').Value = '3.590.91'
$ws.Range('E26').Value = '  -3.85%  '

# Row 27 - PEPE
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000106'
$ws.Range('E27').Value = '  -11.30%  '

# Row 28 - Binance-PegBSC-USD
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.27%  '

# Row 29 - RenderToken
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.29'
$ws.Range('E29').Value = '  -10.26%  '

# Row 30 - PancakeSwap
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.27'
$ws.Range('E30').Value = '  -8.83%  '

# Row 31 - InternetComputer(DFINITY)
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.21'
$ws.Range('E31').Value = '  -12.09%  '

# Row 32 - RenzoRestakedETH
$ws.Range('D32').Value = '3.453.44'
$ws.Range('E32').Value = '  -3.69%  '

# Row 34 - Kaspa
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.146'
$ws.Range('E34').Value = '  -6.98%  '

# Row 35 - EthereumClassic
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '22.87'
$ws.Range('E35').Value = '  -8.56%  '

# Row 36 - Monero
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '173.38'
$ws.Range('E36').Value = '  -1.05%  '

# Row 37 - Fetch.AI
$ws.Range('E37').Value = '  -13.37%  '

# Row 38 - Aptos
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.90'
$ws.Range('E38').Value = '  -10.65%  '

# Row 39 - ImmutableX
$ws.Range('E39').Value = '  -8.59%  '

# Row 40 - NEARProtocol
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.81'
$ws.Range('E40').Value = '  -13.40%  '

# Row 41 - Hedera
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0782'
$ws.Range('E41').Value = '  -8.14%  '

# Row 42 - Mantle
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.817'
$ws.Range('E42').Value = '  -7.00%  '

# Row 43 - OKB
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.59'
$ws.Range('E43').Value = '  -5.30%  '

# Row 44 - FirstDigitalUSD
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  -0.04%  '

# Row 45 - Filecoin
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.46'
$ws.Range('E45').Value = '  -13.88%  '

# Row 47 - EnergySwap
$ws.Range('E47').Value = '  -2.70%  '

# Row 48 - ONDO
$ws.Range('E48').Value = '  -1.35%  '

# Row 49 - Cosmos
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.54'
$ws.Range('E49').Value = '  -8.15%  '

# Row 50 - dogwifhat
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.10'
$ws.Range('E50').Value = '  -16.14%  '

# Row 51 - Maker
$ws.Range('D51').Value = '2.212.34'
$ws.Range('E51').Value = '  -7.52%  '
